$d = $word.ActiveDocument

# Locate the paragraph that holds "Docente(s) Responsável(eis)" (a Heading2)
# and insert a new ListBullet paragraph right after it naming the professor,
# mirroring the target diff.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Docente(s)*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the 'Docente(s) Responsável(eis)' paragraph"
}

$tr = $target.Range

# Position the insertion point one character before the end of the
# paragraph (i.e. just before its trailing paragraph mark) so the new
# paragraph's XML is spliced in without disturbing the existing
# paragraph's own mark/properties.
$insertAt = $d.Range($tr.End - 1, $tr.End - 1)

$xmlSnippet = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>6712818 - Mauricio Lamano Ferreira</w:t></w:r></w:p>'

$insertAt.InsertXML($xmlSnippet)
